$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '66.342.83'
$ws.Range('E2').Value = '  +2.15%  '

Set-TextValue $ws.Range('D3') '3.411.99'
$ws.Range('E3').Value = '  +0.86%  '

Set-TextValue $ws.Range('D4') '1.00'
$ws.Range('E4').Value = '  +0.02%  '

Set-TextValue $ws.Range('D5') '567.51'
$ws.Range('E5').Value = '  +1.37%  '

Set-TextValue $ws.Range('D6') '181.38'
$ws.Range('E6').Value = '  +4.47%  '

Set-TextValue $ws.Range('D7') '0.632'
$ws.Range('E7').Value = '  +1.47%  '

Set-TextValue $ws.Range('D8') '3.403.49'
$ws.Range('E8').Value = '  +0.90%  '

$ws.Range('E9').Value = '  -0.06%  '

$ws.Range('E10').Value = '  +6.87%  '

Set-TextValue $ws.Range('D11') '0.640'
$ws.Range('E11').Value = '  +1.68%  '

Set-TextValue $ws.Range('D12') '54.75'
$ws.Range('E12').Value = '  +1.43%  '

Set-TextValue $ws.Range('D13') '0.0000280'
$ws.Range('E13').Value = '  +0.79%  '

Set-TextValue $ws.Range('D14') '9.34'
$ws.Range('E14').Value = '  +2.74%  '

Set-TextValue $ws.Range('D15') '3.968.75'
$ws.Range('E15').Value = '  +1.39%  '

Set-TextValue $ws.Range('D16') '18.34'
$ws.Range('E16').Value = '  +0.47%  '

Set-TextValue $ws.Range('D17') '3.414.23'
$ws.Range('E17').Value = '  +0.87%  '

$ws.Range('E18').Value = '  +0.66%  '

Set-TextValue $ws.Range('D19') '66.255.03'
$ws.Range('E19').Value = '  +2.16%  '

$ws.Range('E20').Value = '  +1.84%  '

$ws.Range('E21').Value = '  +1.24%  '

Set-TextValue $ws.Range('D22') '467.05'
$ws.Range('E22').Value = '  -0.88%  '

Set-TextValue $ws.Range('D23') '4.99'
$ws.Range('E23').Value = '  +0.30%  '

Set-TextValue $ws.Range('D24') '14.67'
$ws.Range('E24').Value = '  +8.78%  '

Set-TextValue $ws.Range('D25') '4.14'
$ws.Range('E25').Value = '  +0.22%  '

Set-TextValue $ws.Range('D26') '89.64'
$ws.Range('E26').Value = '  +3.40%  '

$ws.Range('E27').Value = '  +1.28%  '

Set-TextValue $ws.Range('D28') '10.83'
$ws.Range('E28').Value = '  +0.48%  '

Set-TextValue $ws.Range('D29') '8.86'
$ws.Range('E29').Value = '  +1.24%  '

Set-TextValue $ws.Range('D30') '31.37'
$ws.Range('E30').Value = '  +2.45%  '

$ws.Range('E31').Value = '  +3.66%  '

Set-TextValue $ws.Range('D32') '11.56'
$ws.Range('E32').Value = '  +0.70%  '

Set-TextValue $ws.Range('D33') '585.88'
$ws.Range('E33').Value = '  +2.55%  '

Set-TextValue $ws.Range('D34') '62.56'
$ws.Range('E34').Value = '  +1.88%  '

$ws.Range('E35').Value = '  +1.23%  '

$ws.Range('E36').Value = '  -0.13%  '

Set-TextValue $ws.Range('D37') '0.145'
$ws.Range('E37').Value = '  +4.26%  '

$ws.Range('E38').Value = '  -1.94%  '

Set-TextValue $ws.Range('D39') '36.44'
$ws.Range('E39').Value = '  +2.43%  '

Set-TextValue $ws.Range('D40') '0.384'
$ws.Range('E40').Value = '  +3.95%  '

$ws.Range('E41').Value = '  +1.92%  '

Set-TextValue $ws.Range('D42') '3.125.45'
$ws.Range('E42').Value = '  +1.28%  '

Set-TextValue $ws.Range('D43') '2.93'
$ws.Range('E43').Value = '  +3.04%  '

Set-TextValue $ws.Range('D44') '0.0425'
$ws.Range('E44').Value = '  +2.43%  '

$ws.Range('B45').Value = 'dogwifhat'
$ws.Range('C45').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue $ws.Range('D45') '2.80'
$ws.Range('E45').Value = '  +19.15%  '

$ws.Range('B46').Value = 'Fetch.AI'
$ws.Range('C46').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue $ws.Range('D46') '2.52'
$ws.Range('E46').Value = '  +2.51%  '

$ws.Range('E47').Value = '  +0.15%  '

Set-TextValue $ws.Range('D48') '3.19'
$ws.Range('E48').Value = '  +2.22%  '

Set-TextValue $ws.Range('D49') '1.00'
$ws.Range('E49').Value = '  +0.24%  '

Set-TextValue $ws.Range('D50') '140.75'
$ws.Range('E50').Value = '  +1.71%  '

Set-TextValue $ws.Range('D51') '8.59'
$ws.Range('E51').Value = '  +3.72%  '
